# Auto-generated: update live market-derived Leve profit columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 933.3333
$ws.Range("I4").Value = 933.3333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 933.3333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -819.3333
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 47620816
$ws.Range("J32").Value = 2063.6667
$ws.Range("L32").Value = 2063.6667
$ws.Range("N32").Value = -2715.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 181818340
$ws.Range("I39").Value = 111111260
$ws.Range("J39").Value = 500000160
$ws.Range("K39").Value = 333333780
$ws.Range("L39").Value = 1500000480
$ws.Range("M39").Value = -333333484
$ws.Range("N39").Value = -1500001072

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 36665.727
$ws.Range("I98").Value = 864.6
$ws.Range("J98").Value = 113382.43
$ws.Range("K98").Value = 864.6
$ws.Range("L98").Value = 113382.43
$ws.Range("M98").Value = 633.4
$ws.Range("N98").Value = -116378.43

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 36665.727
$ws.Range("I122").Value = 864.6
$ws.Range("J122").Value = 113382.43
$ws.Range("K122").Value = 2593.8
$ws.Range("L122").Value = 340147.29
$ws.Range("M122").Value = -143.8000000000002
$ws.Range("N122").Value = -345047.29

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 62528292
$ws.Range("I102").Value = 90911470
$ws.Range("J102").Value = 85286.2
$ws.Range("K102").Value = 90911470
$ws.Range("L102").Value = 85286.2
$ws.Range("M102").Value = -90909848
$ws.Range("N102").Value = -88530.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3240
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3240
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3240
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -3734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2671.577
$ws.Range("I86").Value = 2551.6843
$ws.Range("J86").Value = 2997
$ws.Range("K86").Value = 2551.6843
$ws.Range("L86").Value = 2997
$ws.Range("M86").Value = -1428.6843
$ws.Range("N86").Value = -5243

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2671.577
$ws.Range("I89").Value = 2551.6843
$ws.Range("J89").Value = 2997
$ws.Range("K89").Value = 12758.4215
$ws.Range("L89").Value = 14985
$ws.Range("M89").Value = -7142.4215
$ws.Range("N89").Value = -26217

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2976.6296
$ws.Range("I134").Value = 2255.9565
$ws.Range("K134").Value = 6767.869499999999
$ws.Range("M134").Value = -4232.869499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 38812.375
$ws.Range("J140").Value = 38812.375
$ws.Range("L140").Value = 38812.375
$ws.Range("N140").Value = -49172.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 250001120
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3632

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1618.625
$ws.Range("I35").Value = 1421.2858
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 1421.2858
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -1127.2858
$ws.Range("N35").Value = -3588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2809.5186
$ws.Range("I99").Value = 2554.7
$ws.Range("J99").Value = 2959.4119
$ws.Range("K99").Value = 2554.7
$ws.Range("L99").Value = 2959.4119
$ws.Range("M99").Value = -1056.7
$ws.Range("N99").Value = -5955.4119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2809.5186
$ws.Range("I126").Value = 2554.7
$ws.Range("J126").Value = 2959.4119
$ws.Range("K126").Value = 7664.099999999999
$ws.Range("L126").Value = 8878.235700000001
$ws.Range("M126").Value = -5194.099999999999
$ws.Range("N126").Value = -13818.2357

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5972.5
$ws.Range("J88").Value = 5972.5
$ws.Range("L88").Value = 17917.5
$ws.Range("N88").Value = -18773.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 5972.5
$ws.Range("J91").Value = 5972.5
$ws.Range("L91").Value = 17917.5
$ws.Range("N91").Value = -20881.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1453.5264
$ws.Range("J117").Value = 2622.9
$ws.Range("L117").Value = 7868.700000000001
$ws.Range("N117").Value = -14752.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 104706.1
$ws.Range("I129").Value = 300341.7
$ws.Range("J129").Value = 1740
$ws.Range("K129").Value = 901025.1000000001
$ws.Range("L129").Value = 5220
$ws.Range("M129").Value = -896025.1000000001
$ws.Range("N129").Value = -15220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 159917.11
$ws.Range("I140").Value = 215708.92
$ws.Range("K140").Value = 647126.76
$ws.Range("M140").Value = -641946.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5340.64
$ws.Range("I70").Value = 5476.9414
$ws.Range("J70").Value = 5051
$ws.Range("K70").Value = 5476.9414
$ws.Range("L70").Value = 5051
$ws.Range("M70").Value = -5206.9414
$ws.Range("N70").Value = -5591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5340.64
$ws.Range("I73").Value = 5476.9414
$ws.Range("J73").Value = 5051
$ws.Range("K73").Value = 5476.9414
$ws.Range("L73").Value = 5051
$ws.Range("M73").Value = -4540.9414
$ws.Range("N73").Value = -6923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1036
$ws.Range("I102").Value = 1040
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1040
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 582
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3656.8333
$ws.Range("I32").Value = 3385.2
$ws.Range("J32").Value = 5015
$ws.Range("K32").Value = 3385.2
$ws.Range("L32").Value = 5015
$ws.Range("M32").Value = -3068.2
$ws.Range("N32").Value = -5649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1825
$ws.Range("I100").Value = 1766.6666
$ws.Range("K100").Value = 1766.6666
$ws.Range("M100").Value = -1225.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 399.46667
$ws.Range("I100").Value = 247.2
$ws.Range("J100").Value = 704
$ws.Range("K100").Value = 494.4
$ws.Range("L100").Value = 1408
$ws.Range("M100").Value = 46.60000000000002
$ws.Range("N100").Value = -2490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1505345.6
$ws.Range("I122").Value = 1588920.4
$ws.Range("K122").Value = 4766761.199999999
$ws.Range("M122").Value = -4764311.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2451662
$ws.Range("I126").Value = 2451662
$ws.Range("K126").Value = 7354986
$ws.Range("M126").Value = -7352516
